# Update column F ("dSF") values on the active worksheet to reflect
# repulled/recalculated data, per commit message:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F
$updates = @{
    2  = -2
    3  = -9
    4  = 5
    6  = -5
    8  = 3
    9  = -5
    10 = -5
    12 = -2
    13 = -2
    14 = -5
    15 = -3
    16 = 2
    17 = 4
    18 = 0
    19 = -4
    20 = 2
    21 = 1
    24 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
